$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 724.6667
$ws.Range("I6").Value = 724.6667
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2174.0001
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2062.0001
$ws.Range("N6").ClearContents()
$ws.Range("H17").Value = 3520.8667
$ws.Range("J17").Value = 3520.8667
$ws.Range("L17").Value = 10562.6001
$ws.Range("N17").Value = -10898.6001
$ws.Range("H132").Value = 13992.96
$ws.Range("I132").Value = 2693.9583
$ws.Range("K132").Value = 8081.874899999999
$ws.Range("M132").Value = -5551.874899999999
$ws.Range("H138").Value = 2277.1895
$ws.Range("I138").Value = 1279.9584
$ws.Range("K138").Value = 3839.8752
$ws.Range("M138").Value = 1300.1248
$ws.Range("H140").Value = 69707.375
$ws.Range("J140").Value = 69564.28999999999
$ws.Range("L140").Value = 69564.28999999999
$ws.Range("N140").Value = -79924.28999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3872.4285
$ws.Range("I45").Value = 1112
$ws.Range("J45").Value = 4332.5
$ws.Range("K45").Value = 1112
$ws.Range("L45").Value = 4332.5
$ws.Range("M45").Value = -735
$ws.Range("N45").Value = -5086.5
$ws.Range("H63").Value = 3171.8572
$ws.Range("I63").Value = 1841.8
$ws.Range("K63").Value = 1841.8
$ws.Range("M63").Value = -1155.8
$ws.Range("H66").Value = 3171.8572
$ws.Range("I66").Value = 1841.8
$ws.Range("K66").Value = 9209
$ws.Range("M66").Value = -5777
$ws.Range("H97").Value = 291.51724
$ws.Range("I97").Value = 310
$ws.Range("K97").Value = 310
$ws.Range("M97").Value = 186
$ws.Range("H130").Value = 47372
$ws.Range("J130").Value = 47372
$ws.Range("L130").Value = 47372
$ws.Range("N130").Value = -57412

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3678.2144
$ws.Range("J20").Value = 4162.5
$ws.Range("L20").Value = 4162.5
$ws.Range("N20").Value = -4656.5
$ws.Range("H86").Value = 2889.55
$ws.Range("I86").Value = 2547.1875
$ws.Range("K86").Value = 2547.1875
$ws.Range("M86").Value = -1424.1875
$ws.Range("H89").Value = 2889.55
$ws.Range("I89").Value = 2547.1875
$ws.Range("K89").Value = 12735.9375
$ws.Range("M89").Value = -7119.9375
$ws.Range("H107").Value = 864.375
$ws.Range("I107").Value = 865.4
$ws.Range("J107").Value = 862.6667
$ws.Range("K107").Value = 865.4
$ws.Range("L107").Value = 862.6667
$ws.Range("M107").Value = 1054.6
$ws.Range("N107").Value = -4702.6667
$ws.Range("H131").Value = 48900
$ws.Range("J131").Value = 48900
$ws.Range("L131").Value = 48900
$ws.Range("N131").Value = -58980
$ws.Range("H134").Value = 1678.1765
$ws.Range("I134").Value = 1220.5625
$ws.Range("K134").Value = 3661.6875
$ws.Range("M134").Value = -1126.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2454.1
$ws.Range("I16").Value = 2380.125
$ws.Range("K16").Value = 2380.125
$ws.Range("M16").Value = -2093.125
$ws.Range("H31").Value = 6467.4707
$ws.Range("I31").Value = 2940.889
$ws.Range("K31").Value = 2940.889
$ws.Range("M31").Value = -2645.889
$ws.Range("H34").Value = 6467.4707
$ws.Range("I34").Value = 2940.889
$ws.Range("K34").Value = 2940.889
$ws.Range("M34").Value = -2738.889
$ws.Range("H99").Value = 5494.6787
$ws.Range("I99").Value = 3373.6
$ws.Range("J99").Value = 7942.077
$ws.Range("K99").Value = 3373.6
$ws.Range("L99").Value = 7942.077
$ws.Range("M99").Value = -1875.6
$ws.Range("N99").Value = -10938.077
$ws.Range("H113").Value = 2454.1
$ws.Range("I113").Value = 2380.125
$ws.Range("K113").Value = 2380.125
$ws.Range("M113").Value = -210.125
$ws.Range("H126").Value = 5494.6787
$ws.Range("I126").Value = 3373.6
$ws.Range("J126").Value = 7942.077
$ws.Range("K126").Value = 10120.8
$ws.Range("L126").Value = 23826.231
$ws.Range("M126").Value = -7650.799999999999
$ws.Range("N126").Value = -28766.231
$ws.Range("H134").Value = 1927.3684
$ws.Range("I134").Value = 1927.3684
$ws.Range("K134").Value = 5782.1052
$ws.Range("M134").Value = -3247.1052
$ws.Range("H141").Value = 68361.266
$ws.Range("J141").Value = 73956
$ws.Range("L141").Value = 73956
$ws.Range("N141").Value = -84316

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1007.7143
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 6000
$ws.Range("N5").Value = -6224
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 2
$ws.Range("K16").Value = 6
$ws.Range("M16").Value = 167
$ws.Range("H42").Value = 6427.4287
$ws.Range("J42").Value = 7248.6665
$ws.Range("L42").Value = 21745.9995
$ws.Range("N42").Value = -22813.9995
$ws.Range("H135").Value = 1007.7143
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1981.5
$ws.Range("I113").Value = 1547.25
$ws.Range("K113").Value = 1547.25
$ws.Range("M113").Value = 622.75
$ws.Range("H126").Value = 5296.4375
$ws.Range("I126").Value = 2374.1667
$ws.Range("K126").Value = 7122.500100000001
$ws.Range("M126").Value = -4652.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6555.4116
$ws.Range("I7").Value = 3437
$ws.Range("K7").Value = 3437
$ws.Range("M7").Value = -3325
$ws.Range("H55").Value = 431.125
$ws.Range("I55").Value = 124
$ws.Range("J55").Value = 533.5
$ws.Range("K55").Value = 124
$ws.Range("L55").Value = 533.5
$ws.Range("M55").Value = 49
$ws.Range("N55").Value = -879.5
$ws.Range("H61").Value = 4775.357
$ws.Range("I61").Value = 4097.778
$ws.Range("J61").Value = 5995
$ws.Range("K61").Value = 4097.778
$ws.Range("L61").Value = 5995
$ws.Range("M61").Value = -3895.778
$ws.Range("N61").Value = -6399
$ws.Range("H100").Value = 7980.778
$ws.Range("I100").Value = 2801.25
$ws.Range("K100").Value = 2801.25
$ws.Range("M100").Value = -2260.25
$ws.Range("H113").Value = 4775.357
$ws.Range("I113").Value = 4097.778
$ws.Range("J113").Value = 5995
$ws.Range("K113").Value = 4097.778
$ws.Range("L113").Value = 5995
$ws.Range("M113").Value = -1927.778
$ws.Range("N113").Value = -10335
$ws.Range("H122").Value = 7832.7915
$ws.Range("I122").Value = 4590.25
$ws.Range("J122").Value = 14317.875
$ws.Range("K122").Value = 13770.75
$ws.Range("L122").Value = 42953.625
$ws.Range("M122").Value = -11320.75
$ws.Range("N122").Value = -47853.625
$ws.Range("H126").Value = 6555.4116
$ws.Range("I126").Value = 3437
$ws.Range("K126").Value = 10311
$ws.Range("M126").Value = -7841

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 29900
$ws.Range("J41").Value = 29900
$ws.Range("L41").Value = 29900
$ws.Range("N41").Value = -30680
$ws.Range("H100").Value = 822956.7
$ws.Range("I100").Value = 1015613.25
$ws.Range("K100").Value = 2031226.5
$ws.Range("M100").Value = -2030685.5
$ws.Range("H126").Value = 2728.5
$ws.Range("I126").Value = 1972
$ws.Range("J126").Value = 4998
$ws.Range("K126").Value = 5916
$ws.Range("L126").Value = 14994
$ws.Range("M126").Value = -3446
$ws.Range("N126").Value = -19934
